$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "S-1"
$ws.Range("A3").Value = "S-2"
$ws.Range("A4").Value = "S-3"
$ws.Range("A5").Value = "S-4"
$ws.Range("A6").Value = "S-5"
$ws.Range("A7").Value = "S-6"
$ws.Range("A8").Value = "S-7"
$ws.Range("A9").Value = "S-8"

$ws.Range("A21").Value = "S-1_2"
$ws.Range("A22").Value = "S-2_2"
$ws.Range("A23").Value = "S-3_2"
$ws.Range("A24").Value = "S-4_2"
$ws.Range("A25").Value = "S-5_2"
$ws.Range("A26").Value = "S-6_2"
$ws.Range("A27").Value = "S-7_2"
$ws.Range("A28").Value = "S-8_2"

$ws.Range("C32").Select()
